$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.256.63"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.276.73"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.02"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +1.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.618"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.69%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.596"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0895"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.962"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.624.14"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.279.15"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.581.45"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.38"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.15"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.59"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.34"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.80"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +11.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.92"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.56"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.32"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0843"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.129"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.45"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.28%  "
$ws.Range("E38").Value = "  -2.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.62"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.68"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.84"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.59%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "68.18"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.223"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.81"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.698.05"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.30"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.05"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.97%  "
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.09"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.70%  "
